$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table with
# the latest scraped figures. A handful of Price values read as plain decimals
# (e.g. "128.10", "0.07640") - left as a normal assignment Excel would silently
# coerce them to numbers and drop the significant trailing zero, so those are
# entered with a leading apostrophe to force literal Text, then the cell style
# is reset to Normal so no stray number-format is left behind.

$ws.Range('D2').Value = '28.480.52'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.826.49'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'316.43"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').Value = "'0.5164"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.42%  '
$ws.Range('D8').Value = "'0.3869"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('E9').Value = '  +8.33%  '
$ws.Range('D10').Value = "'1.122"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('D11').Value = "'41.92"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = "'6.394"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.86%  '
$ws.Range('D13').Value = "'21.18"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = "'1.003"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = "'7.496"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('D16').Value = '1.830.12'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = "'94.04"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').Value = "'0.00001126"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.75%  '
$ws.Range('D19').Value = "'0.06637"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = "'17.81"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = "'6.067"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').Value = '28.522.49'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = "'11.45"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.10%  '
$ws.Range('D25').Value = "'2.248"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').Value = "'21.13"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.65%  '
$ws.Range('D27').Value = "'159.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.45%  '
$ws.Range('D28').Value = '2.038.45'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = "'2.420"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').Value = "'125.92"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = "'0.1094"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').Value = "'1.099"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.07640"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.95%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = "'5.735"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('D35').Value = "'3.682"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = "'0.2237"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = "'0.02374"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('D38').Value = "'5.289"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.87%  '
$ws.Range('D39').Value = "'11.99"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.06%  '
$ws.Range('D40').Value = "'8.767"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').Value = "'0.6410"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.63%  '
$ws.Range('D42').Value = "'1.191"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('D43').Value = "'1.399"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').Value = "'13.71"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('D45').Value = "'0.6223"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.46%  '
$ws.Range('D46').Value = "'3.813"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.52%  '
$ws.Range('D47').Value = "'128.10"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('D48').Value = "'2.002"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = "'0.06978"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').Value = "'74.30"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.48%  '
